# B1--and-B2-PowerPoint.pptx
# Commit: Thu, May 07, 2020  8:13:46 PM
#
# The table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") had a new
# Table Quick Style applied to it from PowerPoint's built-in table style
# gallery (Table Tools > Design > Table Styles). This swaps the table's
# <a:tableStyleId> from the deck's custom "Table_0" style
# ({D2BBBFAF-89BB-4984-BAA1-479FAD7AE0A7}) to the built-in
# "Medium Style 2 - Accent 1" style ({207E64B3-600C-4452-9997-459168A829FD}).

$p = $ppt.ActivePresentation

# Locate the slide / shape holding the table (slide 5, shape 2 in this deck)
# by walking the slides instead of hard-coding indices, so the script keeps
# working even if shape ordering shifts slightly.
$targetTable = $null
for ($si = 1; $si -le $p.Slides.Count -and $targetTable -eq $null; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $targetTable = $shape.Table
            break
        }
    }
}

if ($targetTable -ne $null) {
    # Apply the built-in "Medium Style 2 - Accent 1" table style (brace-GUID
    # StyleId), matching what PowerPoint writes when that quick style is
    # picked from the Table Styles gallery.
    $targetTable.ApplyStyle("{207E64B3-600C-4452-9997-459168A829FD}")
}
